$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serial numbers (Excel 1900 date system) corresponding to Dec-31 of each
# year from 1987 through 2024, replacing the old "YYYYQ4" text labels in A2:A39.
$serials = @(
    32142, 32508, 32873, 33238, 33603, 33969, 34334, 34699, 35064, 35430,
    35795, 36160, 36525, 36891, 37256, 37621, 37986, 38352, 38717, 39082,
    39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735,
    43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657
)

for ($i = 0; $i -lt $serials.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $serials[$i]
}

$ws.Range("A2:A39").NumberFormat = "YYYY-MM-DD HH:MM:SS"
